# Add the "03.03" row of data to the summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 - mirror the style used for previous date cells (column A)
$ws.Range("A11").Style = $ws.Range("A10").Style
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "03.03"

$ws.Range("B11").Value = 1034
$ws.Range("C11").Value = 229
$ws.Range("D11").Value = 1000
$ws.Range("E11").Value = 2263
$ws.Range("F11").Value = 160
$ws.Range("G11").Value = 79
$ws.Range("H11").Value = 2502
$ws.Range("I11").Value = 25856
